$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose text looks like a plain decimal number need to be
# forced to Text format first, otherwise Excel COM auto-converts the assigned
# string into a numeric value (losing the intended text representation).
$numericLookingDCells = @(
    "D5", "D8", "D9", "D10", "D11", "D15", "D18", "D19", "D21", "D22",
    "D24", "D25", "D27", "D30", "D31", "D32", "D33", "D36", "D41", "D42",
    "D43", "D45", "D47", "D48"
)
foreach ($addr in $numericLookingDCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '34.119.07'
$ws.Range("E2").Value = '  -0.89%  '

$ws.Range("D3").Value = '1.787.91'
$ws.Range("E3").Value = '  -2.82%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = '224.77'
$ws.Range("E5").Value = '  -0.23%  '

$ws.Range("E6").Value = '  -1.38%  '

$ws.Range("E7").Value = '  -0.12%  '

$ws.Range("D8").Value = '32.94'
$ws.Range("E8").Value = '  +2.92%  '

$ws.Range("D9").Value = '0.286'
$ws.Range("E9").Value = '  -2.35%  '

$ws.Range("D10").Value = '0.0708'
$ws.Range("E10").Value = '  -0.47%  '

$ws.Range("D11").Value = '0.0930'
$ws.Range("E11").Value = '  -0.56%  '

$ws.Range("D12").Value = '2.045.42'
$ws.Range("E12").Value = '  -2.94%  '

$ws.Range("D13").Value = '1.799.83'
$ws.Range("E13").Value = '  -2.26%  '

$ws.Range("E14").Value = '  +0.14%  '

$ws.Range("D15").Value = '0.624'
$ws.Range("E15").Value = '  -3.74%  '

$ws.Range("D16").Value = '34.070.26'
$ws.Range("E16").Value = '  -1.15%  '

$ws.Range("E17").Value = '  -4.73%  '

$ws.Range("D18").Value = '67.91'
$ws.Range("E18").Value = '  -2.73%  '

$ws.Range("D19").Value = '245.95'
$ws.Range("E19").Value = '  -2.19%  '

$ws.Range("D20").Value = '0.0₃0790'
$ws.Range("E20").Value = '  -1.22%  '

$ws.Range("D21").Value = '0.999'
$ws.Range("E21").Value = '  +0.02%  '

$ws.Range("D22").Value = '10.83'
$ws.Range("E22").Value = '  -4.06%  '

$ws.Range("E23").Value = '  -4.24%  '

$ws.Range("D24").Value = '2.11'
$ws.Range("E24").Value = '  -2.63%  '

$ws.Range("D25").Value = '160.61'
$ws.Range("E25").Value = '  -0.64%  '

$ws.Range("E26").Value = '  -2.18%  '

$ws.Range("D27").Value = '7.08'
$ws.Range("E27").Value = '  -2.51%  '

$ws.Range("E28").Value = '  -2.69%  '

$ws.Range("E29").Value = '  -0.10%  '

$ws.Range("D30").Value = '0.0515'
$ws.Range("E30").Value = '  -3.73%  '

$ws.Range("D31").Value = '1.22'
$ws.Range("E31").Value = '  +0.27%  '

$ws.Range("D32").Value = '3.67'
$ws.Range("E32").Value = '  -3.85%  '

$ws.Range("D33").Value = '3.52'
$ws.Range("E33").Value = '  -2.92%  '

$ws.Range("E34").Value = '  -5.96%  '

$ws.Range("D35").Value = '1.398.41'
$ws.Range("E35").Value = '  -4.07%  '

$ws.Range("D36").Value = '0.643'
$ws.Range("E36").Value = '  -0.64%  '

$ws.Range("E37").Value = '  -1.67%  '

$ws.Range("E38").Value = '  -3.74%  '

$ws.Range("E39").Value = '  +3.23%  '

$ws.Range("E40").Value = '  -0.48%  '

$ws.Range("D41").Value = '0.917'
$ws.Range("E41").Value = '  -5.44%  '

$ws.Range("D42").Value = '2.70'
$ws.Range("E42").Value = '  -3.28%  '

$ws.Range("D43").Value = '78.17'
$ws.Range("E43").Value = '  -5.09%  '

$ws.Range("D44").Value = '0.0₆0145'
$ws.Range("E44").Value = '  +15.91%  '

$ws.Range("D45").Value = '1.07'
$ws.Range("E45").Value = '  +1.18%  '

$ws.Range("E46").Value = '  +0.06%  '

$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '107.93'
$ws.Range("E47").Value = '  +0.98%  '

$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").Value = '12.44'
$ws.Range("E48").Value = '  +1.71%  '

$ws.Range("E49").Value = '  -4.51%  '

$ws.Range("D50").Value = '1.944.72'
$ws.Range("E50").Value = '  -3.17%  '

$ws.Range("E51").Value = '  -0.51%  '

# Restore the default cell style for the cells we temporarily reformatted as
# Text, so the resulting workbook does not carry stray formatting changes.
foreach ($addr in $numericLookingDCells) {
    $ws.Range($addr).Style = "Normal"
}
